$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "fixed" marker cell (red-filled H21) and its note (I21 = "исправленно")
$ws.Range("H21").Clear()
$ws.Range("I21").Clear()

# Remove the helper-comment cells in column K (notes explaining why some rows need "+")
$ws.Range("K24").Clear()
$ws.Range("K25").Clear()
$ws.Range("K28").Clear()

# Flip a handful of +/- answers from "+" to "-"
$ws.Range("I24").Value = "-"
$ws.Range("I25").Value = "-"
$ws.Range("H27").Value = "-"
$ws.Range("H28").Value = "-"

# Replace two "+" cells with a new "+/-" value, stored as text (number format "@"),
# matching the new dedicated style used for these cells
$ws.Range("G31").Value = "+/-"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G33").Value = "+/-"
$ws.Range("G33").NumberFormat = "@"
